# chore: prepare Win7 build - rounding/sorting/print fixes
# Converts the monthly report sheet into a daily report sheet:
#  - rename sheet, retitle header, shrink the "统计摘要" font, drop the
#    old bold/white-on-blue "收费项目明细" banner + header styling,
#    insert a new per-day detail table, renumber the remaining rows,
#    widen columns C-F and re-merge the title.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($ws, $row, $col, [string]$text) {
    # Force plain text so Excel's smart-parsing never turns date-like or
    # percent-like strings into numbers/dates with their own number format.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-Num($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

# ---- rename sheet ----
$ws.Name = "日度统计_2025-12"

# ---- wipe existing content/formatting but keep the sheet around ----
$ws.Cells.Clear()

# ---- row 1: title ----
Set-Text $ws 1 1 "2025-12 日度收费统计报表"
$ws.Cells.Item(1,1).Font.Bold = $true
$ws.Cells.Item(1,1).Font.Size = 16
$ws.Cells.Item(1,1).HorizontalAlignment = -4108

# ---- row 2: summary banner ----
Set-Text $ws 2 1 "统计摘要"
$ws.Cells.Item(2,1).Font.Bold = $true
$ws.Cells.Item(2,1).Font.Size = 12

# ---- rows 3-9: summary key/value pairs ----
Set-Text $ws 3 1 "总账单数"
Set-Num  $ws 3 2 45

Set-Text $ws 4 1 "已缴费数"
Set-Num  $ws 4 2 1

Set-Text $ws 5 1 "未缴费数"
Set-Num  $ws 5 2 44

Set-Text $ws 6 1 "总金额"
Set-Text $ws 6 2 "¥107015.40"

Set-Text $ws 7 1 "已缴费金额"
Set-Text $ws 7 2 "¥1369.88"

Set-Text $ws 8 1 "欠费金额"
Set-Text $ws 8 2 "¥105645.52"

Set-Text $ws 9 1 "缴费率"
Set-Text $ws 9 2 "2.2%"

# row 10 intentionally left blank (spacer)

# ---- row 11: daily detail headers ----
Set-Text $ws 11 1 "日期"
Set-Text $ws 11 2 "账单数"
Set-Text $ws 11 3 "日合计(¥)"
Set-Text $ws 11 4 "已缴(¥)"
Set-Text $ws 11 5 "欠费(¥)"

# ---- row 12: daily detail data ----
Set-Text $ws 12 1 "2025-12-16"
Set-Num  $ws 12 2 45
Set-Text $ws 12 3 "¥107015.40"
Set-Text $ws 12 4 "¥1369.88"
Set-Text $ws 12 5 "¥105645.52"

# row 13 intentionally left blank (spacer)

# ---- row 14: fee-item headers ----
Set-Text $ws 14 1 "收费项目"
Set-Text $ws 14 2 "账单数"
Set-Text $ws 14 3 "已缴费数"
Set-Text $ws 14 4 "总金额"
Set-Text $ws 14 5 "已缴金额"
Set-Text $ws 14 6 "欠费金额"

# ---- row 15: fee-item data ----
Set-Text $ws 15 1 "物业费"
Set-Num  $ws 15 2 44
Set-Num  $ws 15 3 0
Set-Text $ws 15 4 "¥106655.40"
Set-Text $ws 15 5 "¥1009.88"
Set-Text $ws 15 6 "¥105645.52"

# ---- row 16: generation timestamp / totals ----
Set-Text $ws 16 1 "生成时间: 2025-12-27 12:09:58"
Set-Num  $ws 16 2 1
Set-Num  $ws 16 3 1
Set-Text $ws 16 4 "¥360.00"
Set-Text $ws 16 5 "¥360.00"
Set-Text $ws 16 6 "¥0.00"

# ---- column widths ----
$ws.Columns.Item(1).ColumnWidth = 18 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 12 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 18 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 18 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 18 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 18 - (5/6)

# ---- merged title ----
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:E1").Merge()

$ws.Range("A1").Select()
